$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 4032.414564323326
$ws.Range("R2").Value = 36291.73107890993
$ws.Range("S2").Value = 0.004461494464728317
$ws.Range("T2").Value = 0.004461494464728317
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 6643.056911897224
$ws.Range("R3").Value = 59787.51220707502
$ws.Range("S3").Value = 0.007349929221942031
$ws.Range("T3").Value = 0.007349929221942032
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 2665.113126663076
$ws.Range("R4").Value = 23986.01813996768
$ws.Range("S4").Value = 0.002948701645828274
$ws.Range("T4").Value = 0.002948701645828275
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 5879.605741311096
$ws.Range("R5").Value = 52916.45167179986
$ws.Range("S5").Value = 0.006505240979369942
$ws.Range("T5").Value = 0.006505240979369943
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 167149.642782898
$ws.Range("R6").Value = 1504346.785046082
$ws.Range("S6").Value = 0.184935649388608
$ws.Range("T6").Value = 0.184935649388608
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 275364.6908316866
$ws.Range("R7").Value = 2478282.217485179
$ws.Range("S7").Value = 0.3046656102268477
$ws.Range("T7").Value = 0.3046656102268477
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 110472.9436896327
$ws.Range("R8").Value = 994256.4932066945
$ws.Range("S8").Value = 0.1222281139281244
$ws.Range("T8").Value = 0.1222281139281244
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 243718.4926518943
$ws.Range("R9").Value = 2193466.433867049
$ws.Range("S9").Value = 0.2696520133467036
$ws.Range("T9").Value = 0.2696520133467037
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 9990.19579822529
$ws.Range("R10").Value = 89911.7621840276
$ws.Range("S10").Value = 0.01105322940991154
$ws.Range("T10").Value = 0.01105322940991154
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 16457.99016692715
$ws.Range("R11").Value = 148121.9115023444
$ws.Range("S11").Value = 0.01820924680709764
$ws.Range("T11").Value = 0.01820924680709765
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 6602.744220633584
$ws.Range("R12").Value = 59424.69798570225
$ws.Range("S12").Value = 0.00730532695050836
$ws.Range("T12").Value = 0.007305326950508362
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 14566.56096120519
$ws.Range("R13").Value = 131099.0486508467
$ws.Range("S13").Value = 0.01611655499747711
$ws.Range("T13").Value = 0.01611655499747711
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 8451.311536473071
$ws.Range("R14").Value = 76061.80382825763
$ws.Range("S14").Value = 0.009350596035751717
$ws.Range("T14").Value = 0.009350596035751717
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 13922.81042075479
$ws.Range("R15").Value = 125305.2937867931
$ws.Range("S15").Value = 0.01540430445203566
$ws.Range("T15").Value = 0.01540430445203566
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 5585.661135303718
$ws.Range("R16").Value = 50270.95021773345
$ws.Range("S16").Value = 0.006180018408198435
$ws.Range("T16").Value = 0.006180018408198436
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 12322.73592876338
$ws.Range("R17").Value = 110904.6233588704
$ws.Range("S17").Value = 0.01363396973686715
$ws.Range("T17").Value = 0.01363396973686716
